$wb = $excel.ActiveWorkbook

# --- Repayment Schedule sheet: insert a new (blank) column before column N ---
# This shifts the existing columns N,O,P to O,P,Q respectively while leaving
# the brand new column N empty, matching the Loan RBI "Variable Instalments"
# column that was added to the repayment schedule.
$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")
$wsSchedule.Columns("N").Insert()

# Give the freshly inserted column roughly the same width as its neighboring
# "In Advance" column (column M), closest achievable to 11.140625 characters.
$wsSchedule.Columns("N").ColumnWidth = 10.33

# --- Sheet/selection activation changes ---
# Make "Repayment Schedule" the active sheet (was "Transactions"), and select
# cell R11 on it (previously A2:XFD11 / A11 was selected).
$wsSchedule.Activate() | Out-Null
$wsSchedule.Range("R11").Select() | Out-Null
